$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: rename the "Diferenciador N" columns (K1:N1) to "Icono N" ---
# (The "Imagen diferenciadora" columns G1:I1 keep their text; the shared-
# string table reshuffles automatically once the old "Diferenciador N"
# strings are no longer referenced anywhere.)
$ws.Range("K1").Value2 = "Icono 1"
$ws.Range("L1").Value2 = "Icono 2"
$ws.Range("M1").Value2 = "Icono 3"
$ws.Range("N1").Value2 = "Icono 4"

# --- New empty, bold-formatted cells in column I, rows 6-9 ---
$iconRange = $ws.Range("I6:I9")
$iconRange.Font.Bold = $true

# --- View state: scroll so column F is the left-most visible column, and
#     move the active selection to K4 ---
$excel.ActiveWindow.ScrollColumn = 6
[void]$ws.Range("K4").Select()
